# Apply the edits from the commit to match the target workbook state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ---
$ws.Range("B10").Value = '5840643 - Luiz Carlos de Queiroz'
$ws.Range("C10").Value = '5840643 - Luiz Carlos de Queiroz'

# --- Row 13 ---
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14 ---
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Introduction to automatic process control. Symbology. Pressure measurement. Temperature measurement. Level measurement. Flow measurement. Controllers. Final control element. Transmission. Systems automation of industrial processes.'
$ws.Range("C14").Value = 'Introduction to automatic process control. Symbology. Pressure measurement. Temperature measurement. Level measurement. Flow measurement. Controllers. Final control element. Transmission. Systems automation of industrial processes.'

# --- Row 15 ---
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2016'
$ws.Range("C15").Value = '01/01/2016'
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16 ---
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1- INTRODUCTION TO AUTOMATIC PROCESS CONTROL. - Introduction. - Open-loop system and closed-loop system. - Symbology. - Feedback control. - Block diagrams. - Transfer function. - Notions of step response of first order process.
2- PRESSURE MEASUREMENT. - Liquid-column methods.  Elastic element methods. - Electrical methods.
3- TEMPERATURE MEASUREMENT. - Thermocouples. - Resistance thermometers. - Filled-system thermometers. - Bimetal thermometers. - Liquid-in-glass thermometers. - Pyrometers.
4- LEVEL MEASUREMENT. - Float-actuated devices. - Pressure devices.
5- FLOW MEASUREMENT. - Orifice meter, Venturi meter, rotameter. - Magnetic flowmeters. - Coriolis mass flowmeters.
6- CONTROLLERS. - On/off control, proportional control, proportional-plus- integral control, proportional-plus-integral-plus-derivative control. - Programmable logic controller.
7- FINAL CONTROL ELEMENT. - Control valves.
8- TRANSMISSION. - Signal transmission.
9- SYSTEMS INDUSTRIAL PROCESS AUTOMATION.'
$ws.Range("C16").Value = '1- INTRODUCTION TO AUTOMATIC PROCESS CONTROL. - Introduction. - Open-loop system and closed-loop system. - Symbology. - Feedback control. - Block diagrams. - Transfer function. - Notions of step response of first order process.
2- PRESSURE MEASUREMENT. - Liquid-column methods.  Elastic element methods. - Electrical methods.
3- TEMPERATURE MEASUREMENT. - Thermocouples. - Resistance thermometers. - Filled-system thermometers. - Bimetal thermometers. - Liquid-in-glass thermometers. - Pyrometers.
4- LEVEL MEASUREMENT. - Float-actuated devices. - Pressure devices.
5- FLOW MEASUREMENT. - Orifice meter, Venturi meter, rotameter. - Magnetic flowmeters. - Coriolis mass flowmeters.
6- CONTROLLERS. - On/off control, proportional control, proportional-plus- integral control, proportional-plus-integral-plus-derivative control. - Programmable logic controller.
7- FINAL CONTROL ELEMENT. - Control valves.
8- TRANSMISSION. - Signal transmission.
9- SYSTEMS INDUSTRIAL PROCESS AUTOMATION.'

# --- Row 17 ---
$ws.Range("B17:C17").Clear()
$ws.Range("A17").Value = 'Avaliação:'
$ws.Rows.Item(17).AutoFit()

# --- Row 18 ---
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5840643 - Luiz Carlos de Queiroz'
$ws.Range("C18").Value = '5840643 - Luiz Carlos de Queiroz'
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19 ---
$ws.Range("A19").Value = 'Critério:'

# --- Row 20 ---
$ws.Range("A20").Value = 'Norma de recuperação:'

# --- Row 21 ---
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22 ---
$ws.Range("B22:C22").Clear()
$ws.Range("A22").Value = 'Requisitos:'
$ws.Rows.Item(22).AutoFit()

# --- Row 23 ---
$ws.Range("A23").Clear()
$ws.Range("B23").Value = 'LOB1006 -  Cálculo IV  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOB1006 -  Cálculo IV  (Requisito fraco)
'
$ws.Rows.Item(23).RowHeight = 30

# --- Row 24 ---
$ws.Range("B24").Value = 'LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOQ4083 -  Fenômenos de Transporte I  (Requisito fraco)
'

# --- Row 25 ---
$ws.Range("B25:C25").Clear()
$ws.Rows.Item(25).AutoFit()

